$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.116.36'
$ws.Range('E2').Value = '  -3.52%  '

$ws.Range('D3').Value = '2.200.76'
$ws.Range('E3').Value = '  -7.08%  '

$c = $ws.Range('D4')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range('E4').Value = '  -0.01%  '

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '295.51'
$c.ClearFormats()
$ws.Range('E5').Value = '  -4.90%  '

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '82.54'
$c.ClearFormats()
$ws.Range('E6').Value = '  -4.05%  '

$ws.Range('E7').Value = '  -3.73%  '

$ws.Range('E8').Value = '  -0.02%  '

$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.465'
$c.ClearFormats()
$ws.Range('E9').Value = '  -5.74%  '

$ws.Range('E10').Value = '  -7.89%  '

$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '28.88'
$c.ClearFormats()
$ws.Range('E11').Value = '  -4.44%  '

$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '47.25'
$c.ClearFormats()
$ws.Range('E12').Value = '  -10.39%  '

$ws.Range('E13').Value = '  -2.14%  '

$ws.Range('D14').Value = '2.543.25'
$ws.Range('E14').Value = '  -6.89%  '

$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '6.23'
$c.ClearFormats()
$ws.Range('E15').Value = '  -4.55%  '

$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '14.06'
$c.ClearFormats()
$ws.Range('E16').Value = '  -6.60%  '

$ws.Range('D17').Value = '2.197.80'
$ws.Range('E17').Value = '  -6.48%  '

$ws.Range('E18').Value = '  -5.88%  '

$ws.Range('D19').Value = '39.000.68'
$ws.Range('E19').Value = '  -3.59%  '

$ws.Range('E20').Value = '  -4.50%  '

$ws.Range('E21').Value = '  -7.32%  '

$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '64.45'
$c.ClearFormats()
$ws.Range('E22').Value = '  -5.82%  '

$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '10.25'
$c.ClearFormats()
$ws.Range('E23').Value = '  -4.45%  '

$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '227.86'
$c.ClearFormats()
$ws.Range('E24').Value = '  -3.07%  '

$ws.Range('E25').Value = '  +0.01%  '

$ws.Range('E26').Value = '  -7.14%  '

$ws.Range('E27').Value = '  -1.62%  '

$ws.Range('E28').Value = '  -5.70%  '

$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '9.07'
$c.ClearFormats()
$ws.Range('E30').Value = '  -1.94%  '

$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '149.10'
$c.ClearFormats()
$ws.Range('E31').Value = '  -3.18%  '

$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '31.65'
$c.ClearFormats()
$ws.Range('E32').Value = '  -7.30%  '

$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range('E33').Value = '  -0.15%  '

$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '4.82'
$c.ClearFormats()
$ws.Range('E34').Value = '  -7.07%  '

$ws.Range('E35').Value = '  -5.02%  '

$ws.Range('E36').Value = '  -3.97%  '

$ws.Range('E37').Value = '  -3.96%  '

$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.0956'
$c.ClearFormats()
$ws.Range('E38').Value = '  -4.71%  '

$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '2.61'
$c.ClearFormats()
$ws.Range('E39').Value = '  -4.86%  '

$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '14.96'
$c.ClearFormats()
$ws.Range('E40').Value = '  -7.37%  '

$ws.Range('E41').Value = '  -5.25%  '

$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '3.64'
$c.ClearFormats()
$ws.Range('E42').Value = '  -4.80%  '

$ws.Range('D43').Value = '1.903.85'

$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '0.0257'
$c.ClearFormats()
$ws.Range('E44').Value = '  -4.01%  '

$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '8.89'
$c.ClearFormats()
$ws.Range('E46').Value = '  -4.94%  '

$ws.Range('E47').Value = '  -2.93%  '

$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '15.77'
$c.ClearFormats()
$ws.Range('E48').Value = '  -11.06%  '

$ws.Range('D49').Value = '2.412.87'
$ws.Range('E49').Value = '  -7.02%  '

$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '70.27'
$c.ClearFormats()
$ws.Range('E50').Value = '  -2.26%  '

$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '86.09'
$c.ClearFormats()
$ws.Range('E51').Value = '  -7.24%  '
